# Fruta / hortaliza, semanal
# A new daily price-report row is inserted at row 57 (Early Majestic /
# Primera, Provincia de Limari, 2021-12-03), pushing every subsequent
# row down by one and extending the used range from T155 to T156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 57; Excel shifts rows 57:155 down
# to 58:156 and carries the row-57 formatting (incl. the date format on
# column D) down onto the newly inserted row.
$ws.Rows("57:57").Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A57").Value = 10
$ws.Range("B57").Value = "Vega Modelo de Temuco"
$ws.Range("C57").Value = "La Araucanía"
$ws.Range("D57").Value = 44533
$ws.Range("E57").Value = 9
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100103
$ws.Range("H57").Value = "Frutos de hueso (carozo)"
$ws.Range("I57").Value = 100103004
$ws.Range("J57").Value = "Durazno"
$ws.Range("K57").Value = "Early Majestic"
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 150
$ws.Range("N57").Value = 24000
$ws.Range("O57").Value = 24000
$ws.Range("P57").Value = 24000
$ws.Range("Q57").Value = "$/bandeja 18 kilos granel"
$ws.Range("R57").Value = "Provincia de Limarí"
$ws.Range("S57").Value = 1333
$ws.Range("T57").Value = 18
